$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 21.35440284743693
$ws.Cells.Item(2, 2).Value = 17.40299466822062
$ws.Cells.Item(2, 3).Value = 25.32300762990845
$ws.Cells.Item(3, 1).Value = 20.46517804228278
$ws.Cells.Item(3, 2).Value = 15.81119853729592
$ws.Cells.Item(3, 3).Value = 25.67267452169189
$ws.Cells.Item(4, 1).Value = 20.5659071169152
$ws.Cells.Item(4, 2).Value = 16.75873551055823
$ws.Cells.Item(4, 3).Value = 25.33821082049982
$ws.Cells.Item(5, 1).Value = 24.85461011005425
$ws.Cells.Item(5, 2).Value = 19.06939777675148
$ws.Cells.Item(5, 3).Value = 30.99478731972179
$ws.Cells.Item(6, 1).Value = 26.19568006624921
$ws.Cells.Item(6, 2).Value = 21.57662977533497
$ws.Cells.Item(6, 3).Value = 30.9694320834139
$ws.Cells.Item(7, 1).Value = 19.27940675326991
$ws.Cells.Item(7, 2).Value = 16.55786402227356
$ws.Cells.Item(7, 3).Value = 22.40840338581645
$ws.Cells.Item(8, 1).Value = 10.16793502270741
$ws.Cells.Item(8, 2).Value = 4.992191239848661
$ws.Cells.Item(8, 3).Value = 15.78123094695072
$ws.Cells.Item(9, 1).Value = 24.45941876439536
$ws.Cells.Item(9, 2).Value = 20.61236004495651
$ws.Cells.Item(9, 3).Value = 29.27056695756565
$ws.Cells.Item(10, 1).Value = 10.12206508708837
$ws.Cells.Item(10, 2).Value = 5.093141499906425
$ws.Cells.Item(10, 3).Value = 16.24157107696975
$ws.Cells.Item(11, 1).Value = 17.9799394048935
$ws.Cells.Item(11, 2).Value = 13.40485233010809
$ws.Cells.Item(11, 3).Value = 22.46111089213312
$ws.Cells.Item(12, 1).Value = 26.94533332817084
$ws.Cells.Item(12, 2).Value = 22.21507864840985
$ws.Cells.Item(12, 3).Value = 31.45720231131714
$ws.Cells.Item(13, 1).Value = 12.64489150394787
$ws.Cells.Item(13, 2).Value = 6.661529136544638
$ws.Cells.Item(13, 3).Value = 20.23858100005091
$ws.Cells.Item(14, 1).Value = 14.72233376916495
$ws.Cells.Item(14, 2).Value = 9.497928009090213
$ws.Cells.Item(14, 3).Value = 20.15579444067186
$ws.Cells.Item(15, 1).Value = 29.03456591231918
$ws.Cells.Item(15, 2).Value = 23.35064153849675
$ws.Cells.Item(15, 3).Value = 35.86897436024218
$ws.Cells.Item(16, 1).Value = 26.24797407641537
$ws.Cells.Item(16, 2).Value = 21.59702931924242
$ws.Cells.Item(16, 3).Value = 31.0327852337292
$ws.Cells.Item(17, 1).Value = 26.97754453567621
$ws.Cells.Item(17, 2).Value = 22.22360870354044
$ws.Cells.Item(17, 3).Value = 31.48400331859357
$ws.Cells.Item(18, 1).Value = 20.59789928049854
$ws.Cells.Item(18, 2).Value = 16.71735006935798
$ws.Cells.Item(18, 3).Value = 25.40335053298117
$ws.Cells.Item(19, 1).Value = 18.64662300382194
$ws.Cells.Item(19, 2).Value = 13.65814444452797
$ws.Cells.Item(19, 3).Value = 24.02813933614454
$ws.Cells.Item(20, 1).Value = 29.58779566904156
$ws.Cells.Item(20, 2).Value = 24.5365450956254
$ws.Cells.Item(20, 3).Value = 34.5246837717676
$ws.Cells.Item(21, 1).Value = 29.00235908482292
$ws.Cells.Item(21, 2).Value = 23.87404234520151
$ws.Cells.Item(21, 3).Value = 34.94931120808728
$ws.Cells.Item(22, 1).Value = 10.60257737430165
$ws.Cells.Item(22, 2).Value = 4.880429883561058
$ws.Cells.Item(22, 3).Value = 17.41572002622064
$ws.Cells.Item(23, 1).Value = 30.78653202353242
$ws.Cells.Item(23, 2).Value = 24.89583517379213
$ws.Cells.Item(23, 3).Value = 37.63742548071873
$ws.Cells.Item(24, 1).Value = 26.24797407641537
$ws.Cells.Item(24, 2).Value = 21.59702931924242
$ws.Cells.Item(24, 3).Value = 31.0327852337292
$ws.Cells.Item(25, 1).Value = 22.46316093330266
$ws.Cells.Item(25, 2).Value = 18.83228079819082
$ws.Cells.Item(25, 3).Value = 26.19568762519166
$ws.Cells.Item(26, 1).Value = 29.34499627619601
$ws.Cells.Item(26, 2).Value = 24.59503716303839
$ws.Cells.Item(26, 3).Value = 34.62307893385925
$ws.Cells.Item(27, 1).Value = 13.00651309338144
$ws.Cells.Item(27, 2).Value = 6.675478500609001
$ws.Cells.Item(27, 3).Value = 20.9646212333834
$ws.Cells.Item(28, 1).Value = 20.39540672253426
$ws.Cells.Item(28, 2).Value = 16.62416845775601
$ws.Cells.Item(28, 3).Value = 24.30745427422584
$ws.Cells.Item(29, 1).Value = 33.12345088628776
$ws.Cells.Item(29, 2).Value = 25.54782596049457
$ws.Cells.Item(29, 3).Value = 39.87908183529395
$ws.Cells.Item(30, 1).Value = 25.66799320215381
$ws.Cells.Item(30, 2).Value = 21.69198451175626
$ws.Cells.Item(30, 3).Value = 30.15833528286886
$ws.Cells.Item(31, 1).Value = 14.38180272393517
$ws.Cells.Item(31, 2).Value = 9.014471796323239
$ws.Cells.Item(31, 3).Value = 19.65954080733745
$ws.Cells.Item(32, 1).Value = 12.07285710960874
$ws.Cells.Item(32, 2).Value = 6.924145520968043
$ws.Cells.Item(32, 3).Value = 18.04051085518306
$ws.Cells.Item(33, 1).Value = 14.75577647820463
$ws.Cells.Item(33, 2).Value = 9.912607647667462
$ws.Cells.Item(33, 3).Value = 20.47474096475314
$ws.Cells.Item(34, 1).Value = 23.25862256463439
$ws.Cells.Item(34, 2).Value = 18.59080453856694
$ws.Cells.Item(34, 3).Value = 28.92472543123027
$ws.Cells.Item(35, 1).Value = 27.99000010796927
$ws.Cells.Item(35, 2).Value = 23.25282439132815
$ws.Cells.Item(35, 3).Value = 32.82710784864798
$ws.Cells.Item(36, 1).Value = 19.23398441712455
$ws.Cells.Item(36, 2).Value = 16.0780844113709
$ws.Cells.Item(36, 3).Value = 22.61114321010599
$ws.Cells.Item(37, 1).Value = 30.76790245543374
$ws.Cells.Item(37, 2).Value = 25.57827995143979
$ws.Cells.Item(37, 3).Value = 36.7197686699047
$ws.Cells.Item(38, 1).Value = 15.15165441290611
$ws.Cells.Item(38, 2).Value = 9.919092188160096
$ws.Cells.Item(38, 3).Value = 21.56595478514924
$ws.Cells.Item(39, 1).Value = 27.40430200393376
$ws.Cells.Item(39, 2).Value = 22.88968686830509
$ws.Cells.Item(39, 3).Value = 31.74665653934783
$ws.Cells.Item(40, 1).Value = 24.81046989252902
$ws.Cells.Item(40, 2).Value = 20.24863346008688
$ws.Cells.Item(40, 3).Value = 29.85225864621885
$ws.Cells.Item(41, 1).Value = 18.01019939613991
$ws.Cells.Item(41, 2).Value = 15.10321019505389
$ws.Cells.Item(41, 3).Value = 20.70120880593804
$ws.Cells.Item(42, 1).Value = 18.22207605481543
$ws.Cells.Item(42, 2).Value = 13.19621269595527
$ws.Cells.Item(42, 3).Value = 23.39604641897973
$ws.Cells.Item(43, 1).Value = 25.66799320215381
$ws.Cells.Item(43, 2).Value = 21.69198451175626
$ws.Cells.Item(43, 3).Value = 30.15833528286886
$ws.Cells.Item(44, 1).Value = 15.80603639507733
$ws.Cells.Item(44, 2).Value = 12.78969308824441
$ws.Cells.Item(44, 3).Value = 18.8987901791203
$ws.Cells.Item(45, 1).Value = 18.86924146006509
$ws.Cells.Item(45, 2).Value = 13.71086760367978
$ws.Cells.Item(45, 3).Value = 24.69780601820212
$ws.Cells.Item(46, 1).Value = 13.98832080062743
$ws.Cells.Item(46, 2).Value = 9.8891746807078
$ws.Cells.Item(46, 3).Value = 18.53846958621791
$ws.Cells.Item(47, 1).Value = 10.69212141709375
$ws.Cells.Item(47, 2).Value = 4.894988708032145
$ws.Cells.Item(47, 3).Value = 18.57720191585278
$ws.Cells.Item(48, 1).Value = 27.13848927895366
$ws.Cells.Item(48, 2).Value = 22.29053005635303
$ws.Cells.Item(48, 3).Value = 31.76097528945665
$ws.Cells.Item(49, 1).Value = 22.80542572449308
$ws.Cells.Item(49, 2).Value = 18.16759686694186
$ws.Cells.Item(49, 3).Value = 27.76653866406352
$ws.Cells.Item(50, 1).Value = 30.77343342930809
$ws.Cells.Item(50, 2).Value = 25.68158874674743
$ws.Cells.Item(50, 3).Value = 36.67926683009082
$ws.Cells.Item(51, 1).Value = 20.50540048967093
$ws.Cells.Item(51, 2).Value = 15.89553849628366
$ws.Cells.Item(51, 3).Value = 25.76663982255887
$ws.Cells.Item(52, 1).Value = 13.62314447925576
$ws.Cells.Item(52, 2).Value = 9.812231006629736
$ws.Cells.Item(52, 3).Value = 17.70007181420092
$ws.Cells.Item(53, 1).Value = 27.57477181131256
$ws.Cells.Item(53, 2).Value = 22.81404009643808
$ws.Cells.Item(53, 3).Value = 32.14326321058359
$ws.Cells.Item(54, 1).Value = 29.02819160186722
$ws.Cells.Item(54, 2).Value = 23.35895602410338
$ws.Cells.Item(54, 3).Value = 35.8746516918341
$ws.Cells.Item(55, 1).Value = 28.94115023989576
$ws.Cells.Item(55, 2).Value = 24.314218892017
$ws.Cells.Item(55, 3).Value = 34.10840275050771
$ws.Cells.Item(56, 1).Value = 14.40810533684264
$ws.Cells.Item(56, 2).Value = 9.058429515973373
$ws.Cells.Item(56, 3).Value = 19.88280472781732
$ws.Cells.Item(57, 1).Value = 12.01756909840416
$ws.Cells.Item(57, 2).Value = 6.60512805544743
$ws.Cells.Item(57, 3).Value = 18.01276797216739
$ws.Cells.Item(58, 1).Value = 30.7410134825391
$ws.Cells.Item(58, 2).Value = 25.3197571103732
$ws.Cells.Item(58, 3).Value = 36.75547520740651
$ws.Cells.Item(59, 1).Value = 17.14318036364662
$ws.Cells.Item(59, 2).Value = 12.61642758959767
$ws.Cells.Item(59, 3).Value = 20.85241626923195
$ws.Cells.Item(60, 1).Value = 30.81039527714425
$ws.Cells.Item(60, 2).Value = 25.112115917388
$ws.Cells.Item(60, 3).Value = 37.48842293111026
$ws.Cells.Item(61, 1).Value = 14.88676128645266
$ws.Cells.Item(61, 2).Value = 9.641954428323732
$ws.Cells.Item(61, 3).Value = 21.08290404360056
$ws.Cells.Item(62, 1).Value = 17.95531811224322
$ws.Cells.Item(62, 2).Value = 13.38666283639152
$ws.Cells.Item(62, 3).Value = 22.33636965776128
$ws.Cells.Item(63, 1).Value = 15.04126008921414
$ws.Cells.Item(63, 2).Value = 9.48685475705855
$ws.Cells.Item(63, 3).Value = 22.30129071795281
$ws.Cells.Item(64, 1).Value = 30.83455451055686
$ws.Cells.Item(64, 2).Value = 25.14039566445994
$ws.Cells.Item(64, 3).Value = 37.53926007847735
$ws.Cells.Item(65, 1).Value = 11.43711189955691
$ws.Cells.Item(65, 2).Value = 6.601256919929904
$ws.Cells.Item(65, 3).Value = 16.50076222766064
$ws.Cells.Item(66, 1).Value = 22.65683687311056
$ws.Cells.Item(66, 2).Value = 17.68063108330326
$ws.Cells.Item(66, 3).Value = 27.57345099009945
$ws.Cells.Item(67, 1).Value = 22.03616481427795
$ws.Cells.Item(67, 2).Value = 17.87391869419686
$ws.Cells.Item(67, 3).Value = 26.58008719190207
$ws.Cells.Item(68, 1).Value = 27.17610390108888
$ws.Cells.Item(68, 2).Value = 22.03592134245795
$ws.Cells.Item(68, 3).Value = 32.19970016935726
$ws.Cells.Item(69, 1).Value = 10.79283631795828
$ws.Cells.Item(69, 2).Value = 5.455960720949792
$ws.Cells.Item(69, 3).Value = 18.91199476908429
$ws.Cells.Item(70, 1).Value = 24.96850649728603
$ws.Cells.Item(70, 2).Value = 20.33182006662399
$ws.Cells.Item(70, 3).Value = 29.81817576862498
$ws.Cells.Item(71, 1).Value = 19.64256785877452
$ws.Cells.Item(71, 2).Value = 16.88583505025246
$ws.Cells.Item(71, 3).Value = 22.64086843200622
$ws.Cells.Item(72, 1).Value = 22.85077740031533
$ws.Cells.Item(72, 2).Value = 19.55563520380617
$ws.Cells.Item(72, 3).Value = 26.1174633631134
